$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 64 (the duplicate/malformed "deepfake" question captured from PDF 3
# of Portuguese). All rows below shift up by one.
$ws.Rows.Item(64).Delete()

# Reset the "respondidas" counter for row 2 back to 0.
$ws.Range("P2").Value = 0
